$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6 (pushes existing row 6.. down by one)
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the January 2025 data
$ws.Range("B6").Value() = 2025
$ws.Range("C6").Value() = "Ene."
$ws.Range("D6").Value() = 360.34199999999998
$ws.Range("E6").Value() = 32087.850999999999
$ws.Range("F6").Value() = 4304.1899999999996
$ws.Range("G6").Value() = 128.77099999999999

# Resize the table / autofilter to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B5:G90"))

# Update the "last updated" note
$ws.Range("B91").Value() = "Actualización: Enero 2025."
